$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) updates ---
$ws.Range("B37").Value = "Stacks"
$ws.Range("B38").Value = "OKB"
$ws.Range("B43").Value = "Bittensor"
$ws.Range("B44").Value = "Mantle"
$ws.Range("B50").Value = "RenderToken"
$ws.Range("B51").Value = "VeChain"

# --- Column C (Link) updates ---
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "59.097.59"
$ws.Range("D3").Value = "2.589.93"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D5").Value = "'527.49"
$ws.Range("D6").Value = "'139.49"
$ws.Range("D8").Value = "'0.564"
$ws.Range("D9").Value = "2.600.77"
$ws.Range("D12").Value = "'0.332"
$ws.Range("D14").Value = "3.047.00"
$ws.Range("D15").Value = "59.024.06"
$ws.Range("D16").Value = "'20.51"
$ws.Range("D18").Value = "2.585.92"
$ws.Range("D19").Value = "'342.73"
$ws.Range("D20").Value = "'4.32"
$ws.Range("D21").Value = "'10.10"
$ws.Range("D22").Value = "'6.42"
$ws.Range("D28").Value = "'7.06"
$ws.Range("D31").Value = "'5.92"
$ws.Range("D32").Value = "'1.61"
$ws.Range("D34").Value = "'149.27"
$ws.Range("D37").Value = "'1.49"
$ws.Range("D38").Value = "'36.82"
$ws.Range("D39").Value = "'0.828"
$ws.Range("D40").Value = "'0.809"
$ws.Range("D41").Value = "'3.53"
$ws.Range("D43").Value = "'271.63"
$ws.Range("D44").Value = "'0.601"
$ws.Range("D45").Value = "'10.74"
$ws.Range("D46").Value = "'0.0953"
$ws.Range("D47").Value = "'0.0515"
$ws.Range("D49").Value = "1.964.32"
$ws.Range("D50").Value = "'4.61"
$ws.Range("D51").Value = "'0.0222"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -2.91%  "
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  -5.04%  "
$ws.Range("E40").Value = "  -7.05%  "
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("E51").Value = "  -0.06%  "
